$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-CellText "D2" '26.871.75'
Set-CellText "E2" '  -1.87%  '
Set-CellText "D3" '1.810.33'
Set-CellText "E3" '  -0.93%  '
Set-CellText "E4" '  +0.13%  '
Set-CellText "D5" '309.50'
Set-CellText "E5" '  -1.19%  '
Set-CellText "E6" '  +0.09%  '
Set-CellText "D7" '0.4644'
Set-CellText "E7" '  +0.49%  '
Set-CellText "D8" '0.3700'
Set-CellText "E8" '  -2.09%  '
Set-CellText "D9" '0.07365'
Set-CellText "E9" '  -0.67%  '
Set-CellText "D10" '0.8753'
Set-CellText "E10" '  -0.28%  '
Set-CellText "D11" '20.45'
Set-CellText "E11" '  -1.63%  '
Set-CellText "B12" 'Polkadot'
Set-CellText "C12" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-CellText "D12" '5.363'
Set-CellText "E12" '  -1.48%  '
Set-CellText "B13" 'WrappedEther'
Set-CellText "C13" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-CellText "D13" '1.731.48'
Set-CellText "E13" '  -5.22%  '
Set-CellText "D14" '6.506'
Set-CellText "E14" '  -3.20%  '
Set-CellText "B15" 'TRON'
Set-CellText "C15" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-CellText "D15" '0.07044'
Set-CellText "E15" '  -0.46%  '
Set-CellText "B16" 'Litecoin'
Set-CellText "C16" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-CellText "D16" '91.59'
Set-CellText "E16" '  -1.66%  '
Set-CellText "D17" '1.002'
Set-CellText "E17" '  +0.13%  '
Set-CellText "D18" '0.000008700'
Set-CellText "E18" '  -1.19%  '
Set-CellText "D19" '1.001'
Set-CellText "E19" '  +0.07%  '
Set-CellText "E20" '  -2.19%  '
Set-CellText "D21" '26.895.87'
Set-CellText "E21" '  -1.79%  '
Set-CellText "D22" '5.306'
Set-CellText "E22" '  -0.55%  '
Set-CellText "D23" '10.54'
Set-CellText "E23" '  -3.84%  '
Set-CellText "D24" '1.978.66'
Set-CellText "E24" '  -3.64%  '
Set-CellText "D25" '1.903'
Set-CellText "E25" '  -2.53%  '
Set-CellText "D26" '151.52'
Set-CellText "E26" '  +0.22%  '
Set-CellText "D27" '18.37'
Set-CellText "E27" '  -1.21%  '
Set-CellText "D28" '2.153'
Set-CellText "E28" '  -5.01%  '
Set-CellText "D29" '5.310'
Set-CellText "E29" '  -0.74%  '
Set-CellText "D30" '115.87'
Set-CellText "E30" '  -1.02%  '
Set-CellText "D31" '0.08903'
Set-CellText "E31" '  -0.52%  '
Set-CellText "D32" '0.7551'
Set-CellText "E32" '  -5.74%  '
Set-CellText "D33" '1.155'
Set-CellText "E33" '  -3.42%  '
Set-CellText "D34" '4.460'
Set-CellText "E34" '  -2.06%  '
Set-CellText "D35" '2.919'
Set-CellText "D36" '1.000'
Set-CellText "D37" '1.099'
Set-CellText "E37" '  +0.00%  '
Set-CellText "D38" '0.01964'
Set-CellText "E38" '  -0.71%  '
Set-CellText "D39" '2.447'
Set-CellText "E39" '  +2.79%  '
Set-CellText "D40" '0.05257'
Set-CellText "E40" '  -0.19%  '
Set-CellText "D41" '2.915'
Set-CellText "E41" '  +0.80%  '
Set-CellText "D42" '0.5317'
Set-CellText "E42" '  -0.50%  '
Set-CellText "D43" '7.171'
Set-CellText "E43" '  -2.67%  '
Set-CellText "D44" '0.1660'
Set-CellText "E44" '  -2.74%  '
Set-CellText "D45" '8.452'
Set-CellText "E45" '  -2.54%  '
Set-CellText "D46" '0.4944'
Set-CellText "E46" '  -3.28%  '
Set-CellText "D47" '10.33'
Set-CellText "E47" '  -2.34%  '
Set-CellText "E48" '  -0.72%  '
Set-CellText "D49" '1.000'
Set-CellText "E49" '  +0.12%  '
Set-CellText "D50" '103.38'
Set-CellText "E50" '  -2.04%  '
Set-CellText "E51" '  -1.38%  '
